$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in "amount filtered" (column H) for the newly-measured rows ---
# Each block of 3 rows shares one filtered-volume reading (per waterbody/date).
$values = @{
    86 = 10;   87 = 10;   88 = 10;
    89 = 2.5;  90 = 2.5;  91 = 2.5;
    92 = 10;   93 = 10;   94 = 10;
    95 = 5;    96 = 5;    97 = 5;
    98 = 5;    99 = 5;    100 = 5;
    101 = 5;   102 = 5;   103 = 5;
    104 = 10;  105 = 10;  106 = 10;
    107 = 10;  108 = 10;  109 = 10;
    110 = 10;  111 = 10;  112 = 10;
    113 = 2.5; 114 = 2.5; 115 = 2.5;
    116 = 2.5;
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 8).Value = $values[$row]
}

# --- Update the sheet view: freeze the header row, zoom to 115%, and move
#     the selection down to where the newly entered data is. ---
$win = $excel.ActiveWindow
[void]$ws.Range("A2").Select()
$win.FreezePanes = $true
$win.Zoom = 115
[void]$ws.Range("H116").Select()
